# Update countries & provincias Spain
# Refreshes the COVID-19 country statistics table on sheet "Pais" to a
# newer snapshot (timestamped 10 de Agosto de 2020 a las 15:08) and fixes
# the ranking order for a few countries whose totals changed enough to
# swap their relative position (Irlanda/Kenia, Corea del Sur/Estado de
# Palestina/Dinamarca, Grecia/Croacia, Timor Oriental/Santa Lucia,
# Montserrat/Islas Malvinas).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp caption
$ws.Range("A1").Value = "Datos actualizados a 10 de Agosto de 2020 a las 15:08"

# Estados Unidos
$ws.Range("B4").Value = 5200394
$ws.Range("C4").Value = 950
$ws.Range("D4").Value = 2664955
$ws.Range("E4").Value = 2369820
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 165619

# India
$ws.Range("B6").Value = 2226229
$ws.Range("C6").Value = 12092
$ws.Range("D6").Value = 1544646
$ws.Range("E6").Value = 636986
$ws.Range("G6").Value = 131
$ws.Range("H6").Value = 44597

# Arabia Saudita
$ws.Range("B16").Value = 289947
$ws.Range("C16").Value = 1257
$ws.Range("D16").Value = 253478
$ws.Range("E16").Value = 33270
$ws.Range("G16").Value = 32
$ws.Range("H16").Value = 3199

# Irak
$ws.Range("B24").Value = 153599
$ws.Range("C24").Value = 3484
$ws.Range("D24").Value = 109790
$ws.Range("E24").Value = 38345
$ws.Range("G24").Value = 72
$ws.Range("H24").Value = 5464

# Filipinas
$ws.Range("B25").Value = 136638
$ws.Range("C25").Value = 6958
$ws.Range("D25").Value = 68159
$ws.Range("E25").Value = 66185
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = 2294

# Bielorrusia
$ws.Range("B42").Value = 68947
$ws.Range("C42").Value = 97
$ws.Range("D42").Value = 64991
$ws.Range("E42").Value = 3367
$ws.Range("G42").Value = 2
$ws.Range("H42").Value = 589

# Paises Bajos
$ws.Range("B45").Value = 59194
$ws.Range("C45").Value = 630

# Row 65 now reports Kenia's updated figures (Kenia moves ahead of Irlanda
# in the country ranking)
$ws.Range("A65").Value = "Kenia"
$ws.Range("B65").Value = 26928
$ws.Range("C65").Value = 492
$ws.Range("D65").Value = 13495
$ws.Range("E65").Value = 13010
$ws.Range("G65").Value = 3
$ws.Range("H65").Value = 423

# Row 66 now reports Irlanda's updated figures
$ws.Range("A66").Value = "Irlanda"
$ws.Range("B66").Value = 26712
$ws.Range("D66").Value = 23364
$ws.Range("E66").Value = 1576
$ws.Range("H66").Value = 1772

# Rows 77-79 rotate: Dinamarca moves ahead of Corea del Sur and Estado de
# Palestina
$ws.Range("A77").Value = "Dinamarca"
$ws.Range("B77").Value = 14815
$ws.Range("C77").Value = 373
$ws.Range("D77").Value = 12925
$ws.Range("E77").Value = 1270
$ws.Range("G77").Value = 3
$ws.Range("H77").Value = 620

$ws.Range("A78").Value = "Corea del Sur"
$ws.Range("B78").Value = 14626
$ws.Range("C78").Value = 28
$ws.Range("D78").Value = 13658
$ws.Range("E78").Value = 663
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 305

$ws.Range("A79").Value = "Estado de Palestina"
$ws.Range("B79").Value = 14510
$ws.Range("C79").Value = 302
$ws.Range("D79").Value = 8045
$ws.Range("E79").Value = 6365
$ws.Range("G79").Value = 3
$ws.Range("H79").Value = 100

# Tayikistan
$ws.Range("B93").Value = 7785
$ws.Range("C93").Value = 40
$ws.Range("D93").Value = 6573
$ws.Range("E93").Value = 1150

# Rows 101-102 swap: Croacia moves ahead of Grecia
$ws.Range("A101").Value = "Croacia"
$ws.Range("B101").Value = 5649
$ws.Range("C101").Value = 45
$ws.Range("D101").Value = 4906
$ws.Range("E101").Value = 585
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 158

$ws.Range("A102").Value = "Grecia"
$ws.Range("B102").Value = 5623
$ws.Range("D102").Value = 3804
$ws.Range("E102").Value = 1607
$ws.Range("H102").Value = 212

# Malta
$ws.Range("B151").Value = 1112
$ws.Range("C151").Value = 23
$ws.Range("D151").Value = 688
$ws.Range("E151").Value = 415

# Rows 202-203 swap: Santa Lucia moves ahead of Timor Oriental (figures
# for both are identical, only the name/order changes)
$ws.Range("A202").Value = "Santa Lucia"

$ws.Range("A203").Value = "Timor Oriental"

# Rows 213-214 swap: Islas Malvinas moves ahead of Montserrat
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
